# Updated cryptos list on Fri Jul 14 19:10:29 UTC 2023 with GitHub Actions
# Applies the per-cell value updates described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.230.01"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.920.37"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("D5").Value = "'246.07"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("E6").Value = "  -12.39%  "
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "'0.3215"
$ws.Range("E8").Value = "  -5.32%  "
$ws.Range("D9").Value = "'26.03"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").Value = "'0.06789"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.7902"
$ws.Range("E11").Value = "  -6.63%  "
$ws.Range("D12").Value = "'0.07967"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "1.931.21"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").Value = "'5.378"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "'93.68"
$ws.Range("E15").Value = "  -8.14%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "'259.69"
$ws.Range("E16").Value = "  -5.78%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.246.91"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'14.26"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "'5.841"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").Value = "'0.000007755"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").Value = "2.189.81"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'1.008"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "'1.010"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "'6.788"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'9.559"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "'158.83"
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("D27").Value = "'18.72"
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("E28").Value = "  -15.71%  "
$ws.Range("D29").Value = "'2.211"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'1.360"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'1.550"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'4.384"
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").Value = "'4.174"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").Value = "'0.05043"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").Value = "'1.188"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'0.7457"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'0.01913"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").Value = "'2.768"
$ws.Range("E39").Value = "  -4.61%  "
$ws.Range("D40").Value = "'79.59"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'6.499"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "'2.022"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").Value = "'0.4385"
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("D44").Value = "'1.007"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "'0.8375"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "'101.29"
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("D47").Value = "'9.672"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("D48").Value = "'7.195"
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("D49").Value = "'35.60"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05945"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.471"
$ws.Range("E51").Value = "  +2.52%  "
